$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old"/"_new" suffixed header columns to "_FV2404"/"_FV2410"
# (the "diff" header in column K/11 stays unchanged).
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"

$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# Wrap the used range in a table ("Table1") so the headers become filterable
# column headers, matching the new xl/tables/table1.xml part.
$tableRange = $ws.Range("A1:U94")
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1, top-left of the scrolling pane
# is A2) so it matches the new frozen-pane sheet view.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
